$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 15
$ws.Range("A15").Value = 12.16
$ws.Range("B15").Value = "陳毅芸"
$ws.Range("C15").Value = "清洗資料、存放資料"
$ws.Range("D15").Value = "code完成"
$ws.Range("E15").Value = "匯入VIX month data.csv檔"
$ws.Range("F15").Value = "匯入爬蟲抓取的標題"

# Add new row 16
$ws.Range("A16").Value = 12.23
$ws.Range("B16").Value = "陳毅芸"
$ws.Range("C16").Value = "清洗資料、存放資料"
$ws.Range("D16").Value = "code完成"
$ws.Range("E16").Value = "將爬蟲標題轉成csv匯入mongoDB"
$ws.Range("F16").Value = "完成書面報告及影片"

# Update the active selection to match the saved workbook view state
$ws.Range("C19").Select()
